# Deploy the implementation guide.
$wb = $excel.ActiveWorkbook

# Rename the "Include from Ferlab.bio CodeS" sheet to "Include #0"
$includeSheet = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$includeSheet.Name = "Include #0"

# Update the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")

# Update Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2024-10-02T15:04:17+00:00"

# Update Contact value (row 10, column B)
$ws.Cells.Item(10, 2).Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row after row 10 (Contact) for "Jurisdiction"
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
# A leading quote forces a text-typed, empty-string cell (matches the
# existing blank shared-string cells used elsewhere, e.g. sheet2!B3)
# rather than a truly-empty cell.
$ws.Cells.Item(11, 2).Value = "'"

# Copy formatting from the row above (Contact row) onto the new row so it
# matches the sheet's established data-row style (border/alignment) instead
# of the default style Insert() leaves behind.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(11, 2).PasteSpecial(-4122) # xlPasteFormats
